$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.796.01"
$ws.Range("E2").Value = "  +7.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.426.45"
$ws.Range("E3").Value = "  +6.35%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.51"
$ws.Range("E5").Value = "  +11.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "319.80"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.13"
$ws.Range("E10").Value = "  +11.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.78"
$ws.Range("E12").Value = "  +7.22%  "
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.03"
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.784.08"
$ws.Range("E16").Value = "  +5.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.419.95"
$ws.Range("E17").Value = "  +6.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.773.96"
$ws.Range("E18").Value = "  +7.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("E20").Value = "  +4.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.22"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("E23").Value = "  +5.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.53"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  +8.15%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.66"
$ws.Range("E27").Value = "  +7.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.38"
$ws.Range("E28").Value = "  +6.00%  "
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.28"
$ws.Range("E30").Value = "  +10.50%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.10"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0965"
$ws.Range("E32").Value = "  +13.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.73"
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.99"
$ws.Range("E34").Value = "  +17.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.123"
$ws.Range("E35").Value = "  +11.03%  "
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.00"
$ws.Range("E37").Value = "  +10.94%  "
$ws.Range("E38").Value = "  +13.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.15"
$ws.Range("E39").Value = "  +15.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0367"
$ws.Range("E40").Value = "  +6.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  +17.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.40"
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("E43").Value = "  +6.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.05"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.42"
$ws.Range("E45").Value = "  +11.80%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.60"
$ws.Range("E47").Value = "  +7.80%  "
$ws.Range("E48").Value = "  +14.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.67"
$ws.Range("E49").Value = "  +18.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.42"
$ws.Range("E50").Value = "  +9.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.00"
$ws.Range("E51").Value = "  +3.70%  "
